$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "Renato"
$ws.Range("B24").Value = "T.I"

# Column C holds numeric-looking salary figures but is stored as TEXT
# throughout the sheet (e.g. C4 = "3000"). Assigning the string directly
# would let Excel auto-coerce it to a number, so copy an existing text
# cell with the same value instead to keep it text without touching styles.
$ws.Range("C4").Copy($ws.Range("C24"))

$ws.Range("D24").Value = 330
$ws.Range("E24").Value = 450
$ws.Range("F24").Value = 240
$ws.Range("G24").Value = 300
$ws.Range("H24").Value = 180

$ws.Range("I24").Value = "22/08/2024"
$ws.Range("J24").Value = "09:43"
